# Updated cryptos list with latest prices and volume percentages
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.607.15"
$ws.Range("E2").Value = "  +0.90%  "

$ws.Range("D3").Value = "1.904.56"
$ws.Range("E3").Value = "  +0.00%  "

$ws.Range("E4").Value = "  +0.04%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "239.12"
$c.Style = "Normal"

$ws.Range("E6").Value = "  -0.08%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4736"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.94%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.2865"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +0.21%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.06655"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.46%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "19.80"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +5.55%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "100.77"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -1.80%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.07808"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +1.28%  "

$ws.Range("D13").Value = "1.902.56"
$ws.Range("E13").Value = "  -0.37%  "

$ws.Range("E14").Value = "  -0.62%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.6789"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +0.60%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "285.08"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +9.61%  "

$ws.Range("D17").Value = "30.630.69"
$ws.Range("E17").Value = "  +0.85%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -0.14%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.000007500"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.03%  "

$ws.Range("D20").Value = "2.163.66"
$ws.Range("E20").Value = "  +0.36%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "12.72"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.13%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.415"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.20%  "

$ws.Range("E23").Value = "  +0.13%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "6.279"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.12%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "9.381"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -0.83%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "167.09"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +1.95%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "19.32"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +2.15%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "2.033"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -1.63%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.380"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.32%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.09944"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -1.74%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "4.513"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -1.76%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "1.516"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +0.50%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "4.268"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +1.32%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.04765"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.16%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.7236"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.93%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.109"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -0.52%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.000"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -0.17%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "2.723"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +0.77%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.01904"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -0.82%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "6.750"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +7.46%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "2.574"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.51%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "74.06"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -0.26%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "1.985"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -0.64%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.8736"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +1.22%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "104.91"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -2.29%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.4272"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +0.44%  "

$ws.Range("E47").Value = "  -0.11%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "991.56"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -0.85%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "7.397"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -1.10%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "9.295"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +4.89%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.1186"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -0.21%  "

